$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete rows 12:42 (shift cells up), leaving A1:A11
$ws.Range("A12:A42").EntireRow.Delete() | Out-Null

# Consolidate each token's scattered property rows (name / type / [abilities] / power-toughness)
# into a single row containing a Python-tuple-style string: ('Name', ['Type', ...extras..., 'P/T'])
$ws.Range("A2").Value  = "('Elemental', ['Token Creature — Elemental', 'Haste', '1/1'])"
$ws.Range("A3").Value  = "('Elf Warrior', ['Token Creature — Elf Warrior', '1/1'])"
$ws.Range("A4").Value  = "('Faerie Rogue', ['Token Creature — Faerie Rogue', 'Flying', '1/1'])"
$ws.Range("A5").Value  = "('Giant Warrior', ['Token Creature — Giant Warrior', 'Haste', '4/4'])"
$ws.Range("A6").Value  = "('Goblin Warrior', ['Token Creature — Goblin Warrior', '1/1'])"
$ws.Range("A7").Value  = "('Kithkin Soldier', ['Token Creature — Kithkin Soldier', '1/1'])"
$ws.Range("A8").Value  = "('Rat', ['Token Creature — Rat', '1/1'])"
$ws.Range("A9").Value  = "('Spider', ['Token Creature — Spider', 'Reach', '1/2'])"
$ws.Range("A10").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A11").Value = "('Wolf', ['Token Creature — Wolf', '2/2'])"
